# Add a new cancer type row ("Stomach and esophageal cancer" / "STES")
# to the bottom of the "Cancers" sheet, mirroring the formatting of the
# existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cancers")

# Find the next empty row right after the current data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 2).Value = "STES"
$ws.Cells.Item($newRow, 1).Value = "Stomach and esophageal cancer"

# Match the row height used by the other data rows (e.g. row 38)
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($lastRow).RowHeight

# Leave the sheet scrolled/selected where the author ended up after typing
$null = $ws.Activate()
$null = $ws.Range("A31").Select()
